$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 09:33"

# --- Row 63/64: Armenia / Argelia swap region ---
$ws.Range("A63").Value = "Armenia"
$ws.Range("B63").Value = 51925
$ws.Range("C63").Value = 543
$ws.Range("D63").Value = 44583
$ws.Range("E63").Value = 6370
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 972

$ws.Range("A64").Value = "Argelia"
$ws.Range("B64").Value = 51847
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 36385
$ws.Range("E64").Value = 13713
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 1749

# --- Row 72: Afganistan (label unchanged, new stats) ---
$ws.Range("B72").Value = 39297
$ws.Range("C72").Value = 7
$ws.Range("D72").Value = 32842
$ws.Range("E72").Value = 4993
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 1462

# --- Rows 77/78/79: Hungria / El Salvador / Dinamarca reshuffle ---
$ws.Range("A77").Value = "Hungria"
$ws.Range("B77").Value = 29717
$ws.Range("C77").Value = 1086
$ws.Range("D77").Value = 6824
$ws.Range("E77").Value = 22081
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 14
$ws.Range("H77").Value = 812

$ws.Range("A78").Value = "El Salvador"
$ws.Range("B78").Value = 29175
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 23930
$ws.Range("E78").Value = 4392
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 853

$ws.Range("A79").Value = "Dinamarca"
$ws.Range("B79").Value = 28932
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 21824
$ws.Range("E79").Value = 6456
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 652

# --- Rows 114/115: Georgia / Mauritania swap ---
$ws.Range("A114").Value = "Georgia"
$ws.Range("B114").Value = 7564
$ws.Range("C114").Value = 471
$ws.Range("D114").Value = 3992
$ws.Range("E114").Value = 3524
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 48

$ws.Range("A115").Value = "Mauritania"
$ws.Range("B115").Value = 7511
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 7163
$ws.Range("E115").Value = 186
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 162

# --- Rows 158/159: Letonia / Polinesia Francesa swap ---
$ws.Range("A158").Value = "Letonia"
$ws.Range("B158").Value = 2019
$ws.Range("C158").Value = 74
$ws.Range("D158").Value = 1307
$ws.Range("E158").Value = 674
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 38

$ws.Range("A159").Value = "Polinesia Francesa"
$ws.Range("B159").Value = 1964
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 1555
$ws.Range("E159").Value = 401
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 8
